$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Neurology" to "Session"
$ws.Name = "Session"

# Remove the last data row (row 86), shrinking the used range to A1:F85
$ws.Rows(86).Delete()
